# Actualización automática 2025-06-05 10:49:05
# Adds a new "CUMPLIMIENTO MENSUAL" worksheet summarizing budget vs. sales
# compliance per product group for HIDALGO HIDALGO PEDRO GUSTAVO.

$wb = $excel.ActiveWorkbook

# --- Create the new worksheet, positioned after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "CUMPLIMIENTO MENSUAL"

# --- Column widths (match "width" stored attr = ColumnWidth + 5/6) ---
$offset = 5/6
$ws.Columns.Item(1).ColumnWidth = 31 - $offset
$ws.Columns.Item(2).ColumnWidth = 22 - $offset
$ws.Columns.Item(3).ColumnWidth = 22 - $offset
$ws.Columns.Item(4).ColumnWidth = 11 - $offset
$ws.Columns.Item(5).ColumnWidth = 22 - $offset
$ws.Columns.Item(6).ColumnWidth = 18 - $offset

# --- Header row (reuse the existing bold+bordered header style via copy) ---
$headerSrc = $wb.Worksheets.Item(1).Range("A1")

$ws.Range("A1").Value = "ASESOR"
$ws.Range("B1").Value = "GRUPO"
$ws.Range("C1").Value = "PRESUPUESTO"
$ws.Range("D1").Value = "VENTA"
$ws.Range("E1").Value = "POR CUMPLIR"
$ws.Range("F1").Value = "CUMPLIMIENTO"

$headerSrc.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# --- Data rows ---
$asesor = "HIDALGO HIDALGO PEDRO GUSTAVO"

$grupos = @(
    @("240X120 PORCELANATO", 782.465010521559, 0, 782.465010521559),
    @("240X80 PORCELANATO", 4168.07156573679, 0, 4168.07156573679),
    @("FREGADEROS DE COCINA", 513.831046659336, 0, 513.831046659336),
    @("GRANITO", 238.32, 0, 238.32),
    @("GRIFERIAS", 106.82, 0, 106.82),
    @("INODOROS", 1800, 0, 1800),
    @("LAVABOS", 625, 0, 625),
    @("LED", 300, 0, 300),
    @("NO RESURTIBLES", 650.25, 0, 650.25),
    @("OTROS", 0, 0, 0),
    @("PANELES DECORATIVOS", 350, 0, 350),
    @("PANELES PU", 230, 0, 230),
    @("PANELES PVC", 483, 0, 483),
    @("PIEDRA SINTERIZADA", 7465, 0, 7465),
    @("PORCELANATO", 29532.44, 0, 29532.44),
    @("PUERTAS DE SEGURIDAD", 342, 0, 342),
    @("SAL SOLUBLE", 2800, 0, 2800)
)

# Existing currency-format cell (style used for PRESUPUESTO/VENTA/POR CUMPLIR columns)
$currencySrc = $wb.Worksheets.Item(2).Range("C2")
# Existing currency-format, right-aligned cell is style idx 6 (not what we want here);
# the CUMPLIMIENTO column needs the plain percent style (numFmtId 10), currently unused
# anywhere in the workbook, so set it directly via NumberFormat.

$row = 2
foreach ($g in $grupos) {
    $ws.Cells.Item($row, 1).Value = $asesor
    $ws.Cells.Item($row, 2).Value = $g[0]
    $ws.Cells.Item($row, 3).Value = $g[1]
    $ws.Cells.Item($row, 4).Value = $g[2]
    $ws.Cells.Item($row, 5).Value = $g[3]
    $ws.Cells.Item($row, 6).Value = 0

    $currencySrc.Copy()
    $ws.Range($ws.Cells.Item($row, 3), $ws.Cells.Item($row, 5)).PasteSpecial(-4122)
    $ws.Cells.Item($row, 6).NumberFormat = "0.00%"

    $row = $row + 1
}

# --- TOTAL row (row 19) ---
$ws.Range("B19").Value = "TOTAL"
$ws.Range("C19").Value = 50387.19762291768
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 50387.19762291768
$ws.Range("F19").Value = 0

$currencySrc.Copy()
$ws.Range("C19:E19").PasteSpecial(-4122)
$ws.Range("F19").NumberFormat = "0.00%"

# TOTAL label style: right-aligned, no border/bold/number-format (new style idx 7)
$ws.Range("B19").HorizontalAlignment = -4152

$null = $ws.Range("A1").Select()

Write-Output "CUMPLIMIENTO MENSUAL sheet created"
